# Insert a new row above the current row 2 ("Groceries"); this shifts
# Groceries/Insurance/Miscellaneous down one row (rows 2-4 -> 3-5) and
# extends the sheet's used range / dimension to A1:C5 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Insert()

# The new row 2 has no formatting yet; copy the date style (numFmt) from
# the row directly below it (row 3, the former "Groceries" row, C3) onto
# the new C2 cell so the new date value renders the same way.
$ws.Cells.Item(3, 3).Copy()
$ws.Cells.Item(2, 3).PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new "Car Insurance" expense row.
$ws.Cells.Item(2, 1).Value = "Car Insurance"
$ws.Cells.Item(2, 2).Value = 1250
$ws.Cells.Item(2, 3).Value = 45837.22928240741
